# Update cryptocurrency price/volume snapshot values (scheduled data refresh),
# including two rows (35/36) whose coin ordering swapped between runs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.983.27"
$ws.Range("E2").Value = "  -1.59%  "
$ws.Range("D3").Value = "1.830.20"
$ws.Range("E3").Value = "  -2.18%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'239.49"
$ws.Range("E5").Value = "  -1.80%  "
$ws.Range("D6").Value = "'0.6855"
$ws.Range("E6").Value = "  -3.28%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.07626"
$ws.Range("E8").Value = "  -3.31%  "
$ws.Range("D9").Value = "'0.3015"
$ws.Range("E9").Value = "  -4.76%  "
$ws.Range("D10").Value = "'23.37"
$ws.Range("E10").Value = "  -5.23%  "
$ws.Range("D11").Value = "'0.07749"
$ws.Range("E11").Value = "  -3.16%  "
$ws.Range("D12").Value = "1.842.22"
$ws.Range("E12").Value = "  -1.88%  "
$ws.Range("D13").Value = "'5.045"
$ws.Range("E13").Value = "  -3.54%  "
$ws.Range("D14").Value = "'90.42"
$ws.Range("E14").Value = "  -3.98%  "
$ws.Range("D15").Value = "'0.6737"
$ws.Range("E15").Value = "  -4.62%  "
$ws.Range("D16").Value = "'6.443"
$ws.Range("E16").Value = "  -1.19%  "
$ws.Range("D17").Value = "'0.000008273"
$ws.Range("E17").Value = "  -1.11%  "
$ws.Range("D18").Value = "28.990.86"
$ws.Range("E18").Value = "  -1.65%  "
$ws.Range("D19").Value = "'242.82"
$ws.Range("E19").Value = "  -5.68%  "
$ws.Range("D20").Value = "2.103.89"
$ws.Range("E20").Value = "  -1.30%  "
$ws.Range("D21").Value = "'12.67"
$ws.Range("E21").Value = "  -4.28%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "'7.442"
$ws.Range("E23").Value = "  -2.57%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").Value = "'0.1469"
$ws.Range("E25").Value = "  -6.03%  "
$ws.Range("D26").Value = "'161.26"
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").Value = "'8.732"
$ws.Range("E27").Value = "  -3.92%  "
$ws.Range("D28").Value = "'18.14"
$ws.Range("E28").Value = "  -4.21%  "
$ws.Range("D29").Value = "'1.531"
$ws.Range("E29").Value = "  +1.91%  "
$ws.Range("E30").Value = "  -3.25%  "
$ws.Range("D31").Value = "'4.163"
$ws.Range("E31").Value = "  -2.38%  "
$ws.Range("D32").Value = "'1.192"
$ws.Range("E32").Value = "  -1.45%  "
$ws.Range("D33").Value = "'0.05119"
$ws.Range("E33").Value = "  -3.84%  "
$ws.Range("D34").Value = "'0.7587"
$ws.Range("E34").Value = "  +0.98%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "'1.818"
$ws.Range("E35").Value = "  -4.47%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'1.150"
$ws.Range("E36").Value = "  -2.28%  "
$ws.Range("D37").Value = "'2.702"
$ws.Range("E37").Value = "  -0.41%  "
$ws.Range("D38").Value = "'0.01834"
$ws.Range("E38").Value = "  -2.71%  "
$ws.Range("D39").Value = "1.219.76"
$ws.Range("E39").Value = "  -4.13%  "
$ws.Range("D40").Value = "'2.710"
$ws.Range("E40").Value = "  -1.64%  "
$ws.Range("D41").Value = "'0.9126"
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("D42").Value = "'108.70"
$ws.Range("D43").Value = "'0.9998"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").Value = "2.003.09"
$ws.Range("E44").Value = "  -1.19%  "
$ws.Range("D45").Value = "'5.416"
$ws.Range("E45").Value = "  -9.55%  "
$ws.Range("D46").Value = "'0.5175"
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("D47").Value = "'9.423"
$ws.Range("E47").Value = "  -1.30%  "
$ws.Range("D48").Value = "'63.42"
$ws.Range("E48").Value = "  -11.62%  "
$ws.Range("E49").Value = "  -8.32%  "
$ws.Range("D50").Value = "'1.728"
$ws.Range("E50").Value = "  -3.81%  "
$ws.Range("D51").Value = "'6.903"
$ws.Range("E51").Value = "  -2.56%  "
